$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45 (pushes old rows 45..169 down to 46..170)
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new data record
$ws.Cells.Item(45,1).Value  = 7
$ws.Cells.Item(45,2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(45,3).Value  = 'Ñuble'
$ws.Cells.Item(45,4).Value  = 44487
$ws.Cells.Item(45,5).Value  = 16
$ws.Cells.Item(45,6).Value  = 100112009
$ws.Cells.Item(45,7).Value  = 'Acelga'
$ws.Cells.Item(45,8).Value  = 'Sin especificar'
$ws.Cells.Item(45,9).Value  = 'Primera'
$ws.Cells.Item(45,10).Value = 120
$ws.Cells.Item(45,11).Value = 350
$ws.Cells.Item(45,12).Value = 400
$ws.Cells.Item(45,13).Value = 375
$ws.Cells.Item(45,14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(45,15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(45,16).Value = 375
$ws.Cells.Item(45,17).Value = 1
$ws.Cells.Item(45,18).Value = 'Hortaliza'
